# Update "想去人数" (interested-count) figures in the 展览 and 全部类型 sheets
# to reflect the newly scraped totals (gh-pages data refresh at 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of row number -> new F-column value for worksheet "展览"
$zhanlanUpdates = @{
    2  = 70
    3  = 384
    4  = 442
    5  = 25
    6  = 12
    7  = 251
    8  = 13867
    9  = 83
    10 = 78
    11 = 5592
    13 = 44
    14 = 35
    16 = 1219
    18 = 161
    19 = 750
    20 = 2903
    22 = 10385
    24 = 24
    25 = 38
    26 = 3701
}

# Map of row number -> new F-column value for worksheet "全部类型"
$quanbuUpdates = @{
    2  = 70
    3  = 384
    5  = 442
    6  = 25
    7  = 12
    8  = 251
    9  = 13867
    10 = 83
    11 = 78
    12 = 5592
    14 = 44
    15 = 35
    17 = 1219
    19 = 161
    20 = 750
    21 = 2903
    24 = 10385
    26 = 24
    27 = 38
    28 = 3701
}

$wsZhanlan = $wb.Worksheets.Item("展览")
foreach ($row in $zhanlanUpdates.Keys) {
    $wsZhanlan.Range("F$row").Value = $zhanlanUpdates[$row]
}

$wsQuanbu = $wb.Worksheets.Item("全部类型")
foreach ($row in $quanbuUpdates.Keys) {
    $wsQuanbu.Range("F$row").Value = $quanbuUpdates[$row]
}
